$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.466.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.579.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.497"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.47%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.17"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0588"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0864"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.806.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.583.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.522"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.490.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0690"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.92%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0471"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.362.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.967"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0167"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.530"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.820"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.970"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("E44").Value = "  +3.65%  "
$ws.Range("E45").Value = "  -1.53%  "
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.717.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₇0998"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0956"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("E51").Value = "  -0.56%  "
